$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3436.375
$ws.Range("J112").Value = 5749
$ws.Range("L112").Value = 17247
$ws.Range("N112").Value = -19463

$ws.Range("H129").Value = 1334.5217
$ws.Range("I129").Value = 1024.25
$ws.Range("K129").Value = 3072.75
$ws.Range("M129").Value = 1927.25

$ws.Range("H132").Value = 5951.7295
$ws.Range("I132").Value = 5462.647
$ws.Range("K132").Value = 16387.941
$ws.Range("M132").Value = -13857.941

$ws.Range("H138").Value = 2886.4666
$ws.Range("J138").Value = 3021.875
$ws.Range("L138").Value = 9065.625
$ws.Range("N138").Value = -19345.625


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H112").Value = 15095.25
$ws.Range("J112").Value = 15095.25
$ws.Range("L112").Value = 15095.25
$ws.Range("N112").Value = -18049.25

$ws.Range("H132").Value = 24901.268
$ws.Range("I132").Value = 25308.182
$ws.Range("J132").Value = 6997
$ws.Range("K132").Value = 75924.546
$ws.Range("L132").Value = 20991
$ws.Range("M132").Value = -73394.546
$ws.Range("N132").Value = -26051


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3913.8
$ws.Range("I99").Value = 3445.4443
$ws.Range("K99").Value = 3445.4443
$ws.Range("M99").Value = -1947.4443

$ws.Range("H105").Value = 2594.6956
$ws.Range("I105").Value = 2325.842
$ws.Range("J105").Value = 3871.75
$ws.Range("K105").Value = 2325.842
$ws.Range("L105").Value = 3871.75
$ws.Range("M105").Value = -578.8420000000001
$ws.Range("N105").Value = -7365.75

$ws.Range("H134").Value = 2714.5676
$ws.Range("I134").Value = 2847.4443
$ws.Range("K134").Value = 8542.332900000001
$ws.Range("M134").Value = -6007.332900000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5288.125
$ws.Range("I12").Value = 5836.3335
$ws.Range("J12").Value = 4959.2
$ws.Range("K12").Value = 5836.3335
$ws.Range("L12").Value = 4959.2
$ws.Range("M12").Value = -5666.3335
$ws.Range("N12").Value = -5299.2

$ws.Range("H31").Value = 7848.7
$ws.Range("I31").Value = 7249.5
$ws.Range("J31").Value = 8248.166999999999
$ws.Range("K31").Value = 7249.5
$ws.Range("L31").Value = 8248.166999999999
$ws.Range("M31").Value = -6954.5
$ws.Range("N31").Value = -8838.166999999999

$ws.Range("H34").Value = 7848.7
$ws.Range("I34").Value = 7249.5
$ws.Range("J34").Value = 8248.166999999999
$ws.Range("K34").Value = 7249.5
$ws.Range("L34").Value = 8248.166999999999
$ws.Range("M34").Value = -7047.5
$ws.Range("N34").Value = -8652.166999999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 190
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 437.25
$ws.Range("J60").Value = 433
$ws.Range("L60").Value = 1299
$ws.Range("N60").Value = -1801

$ws.Range("H61").Value = 623.3333
$ws.Range("I61").Value = 120
$ws.Range("K61").Value = 360
$ws.Range("M61").Value = -145


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 44999
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H70").Value = 4707
$ws.Range("I70").Value = 4433.5
$ws.Range("K70").Value = 4433.5
$ws.Range("M70").Value = -4163.5

$ws.Range("H73").Value = 4707
$ws.Range("I73").Value = 4433.5
$ws.Range("K73").Value = 4433.5
$ws.Range("M73").Value = -3497.5

$ws.Range("H80").Value = 6521
$ws.Range("I80").Value = 2364.75
$ws.Range("J80").Value = 11271
$ws.Range("K80").Value = 2364.75
$ws.Range("L80").Value = 11271
$ws.Range("M80").Value = -1366.75
$ws.Range("N80").Value = -13267

$ws.Range("H81").Value = 44999
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H83").Value = 6521
$ws.Range("I83").Value = 2364.75
$ws.Range("J83").Value = 11271
$ws.Range("K83").Value = 11823.75
$ws.Range("L83").Value = 56355
$ws.Range("M83").Value = -6831.75
$ws.Range("N83").Value = -66339

$ws.Range("H84").Value = 44999
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H132").Value = 278347.1
$ws.Range("I132").Value = 337979.88
$ws.Range("K132").Value = 1013939.64
$ws.Range("M132").Value = -1011409.64


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3898.862
$ws.Range("I40").Value = 3294.1667
$ws.Range("K40").Value = 3294.1667
$ws.Range("M40").Value = -3158.1667

$ws.Range("H81").Value = 61181
$ws.Range("J81").Value = 61181
$ws.Range("L81").Value = 61181
$ws.Range("N81").Value = -63177

$ws.Range("H84").Value = 61181
$ws.Range("J84").Value = 61181
$ws.Range("L84").Value = 183543
$ws.Range("N84").Value = -193527

$ws.Range("H122").Value = 5459.3335
$ws.Range("I122").Value = 5680
$ws.Range("K122").Value = 17040
$ws.Range("M122").Value = -14590

$ws.Range("H136").Value = 6590.364
$ws.Range("I136").Value = 5513.4287
$ws.Range("K136").Value = 16540.2861
$ws.Range("M136").Value = -13990.2861


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 49999
$ws.Range("J59").Value = 49999
$ws.Range("L59").Value = 49999
$ws.Range("N59").Value = -51475

$ws.Range("H93").Value = 99389
$ws.Range("J93").Value = 99389
$ws.Range("L93").Value = 99389
$ws.Range("N93").Value = -104381

$ws.Range("H126").Value = 103609.8
$ws.Range("I126").Value = 114566.445
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 343699.335
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -341229.335
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 94284
$ws.Range("I132").Value = 128655.875
$ws.Range("K132").Value = 385967.625
$ws.Range("M132").Value = -383437.625

$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -100119

$ws.Range("H136").Value = 2953.5334
$ws.Range("J136").Value = 8500
$ws.Range("L136").Value = 25500
$ws.Range("N136").Value = -30600

